# Generate Report for Archive
#
# 1) Status text: "Ready for handoff" -> "In Translation" (Overview sheet,
#    row 2, zh-cn / de-de status columns E2 & F2 both hold this string).
# 2) Narrow the now-shorter "Status" columns:
#      - Overview!E:F  (zh-cn / de-de status columns)
#      - zh-cn!C       (Status column)
#      - de-de!C       (Status column)

$wb = $excel.ActiveWorkbook

# --- 1) Update the status text everywhere it appears (Overview, zh-cn, de-de) ---
$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Cells.Replace("Ready for handoff", "In Translation")
$zhcn.Cells.Replace("Ready for handoff", "In Translation")
$dede.Cells.Replace("Ready for handoff", "In Translation")

# --- 2) Resize the Status columns ---
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
